# Auto-generated edit script: updates Universalis market-price derived
# columns (H-N) across multiple crafting-profession sheets, matching the
# "chore: update Sheets via scheduled runner" price refresh commit.

$wb = $excel.ActiveWorkbook

### Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")

# Row 11
$ws.Range("H11").Value = 1393
$ws.Range("I11").Value = 1393
$ws.Range("K11").Value = 1393
$ws.Range("M11").Value = -1253

# Row 64
$ws.Range("H64").Value = 3840.2942
$ws.Range("I64").Value = 3516.818
$ws.Range("J64").Value = 4433.3335
$ws.Range("K64").Value = 3516.818
$ws.Range("L64").Value = 4433.3335
$ws.Range("M64").Value = -3268.818
$ws.Range("N64").Value = -4929.3335

# Row 67
$ws.Range("H67").Value = 3840.2942
$ws.Range("I67").Value = 3516.818
$ws.Range("J67").Value = 4433.3335
$ws.Range("K67").Value = 3516.818
$ws.Range("L67").Value = 4433.3335
$ws.Range("M67").Value = -2658.818
$ws.Range("N67").Value = -6149.3335

# Row 74
$ws.Range("H74").Value = 5289.25
$ws.Range("I74").Value = 4309.125
$ws.Range("J74").Value = 7249.5
$ws.Range("K74").Value = 4309.125
$ws.Range("L74").Value = 7249.5
$ws.Range("M74").Value = -3373.125
$ws.Range("N74").Value = -9121.5

# Row 76
$ws.Range("H76").Value = 4164.1797
$ws.Range("I76").Value = 3982.862
$ws.Range("J76").Value = 4690
$ws.Range("K76").Value = 3982.862
$ws.Range("L76").Value = 4690
$ws.Range("M76").Value = -3667.862
$ws.Range("N76").Value = -5320

# Row 77
$ws.Range("H77").Value = 5289.25
$ws.Range("I77").Value = 4309.125
$ws.Range("J77").Value = 7249.5
$ws.Range("K77").Value = 21545.625
$ws.Range("L77").Value = 36247.5
$ws.Range("M77").Value = -16865.625
$ws.Range("N77").Value = -45607.5

# Row 79
$ws.Range("H79").Value = 4164.1797
$ws.Range("I79").Value = 3982.862
$ws.Range("J79").Value = 4690
$ws.Range("K79").Value = 3982.862
$ws.Range("L79").Value = 4690
$ws.Range("M79").Value = -2890.862
$ws.Range("N79").Value = -6874

# Row 92
$ws.Range("H92").Value = 298.0909
$ws.Range("I92").Value = 211.28572
$ws.Range("J92").Value = 450
$ws.Range("K92").Value = 211.28572
$ws.Range("L92").Value = 450
$ws.Range("M92").Value = 1036.71428
$ws.Range("N92").Value = -2946

# Row 97
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

# Row 107
$ws.Range("H107").Value = 506.71875
$ws.Range("I107").Value = 411.8846
$ws.Range("J107").Value = 917.6667
$ws.Range("K107").Value = 411.8846
$ws.Range("L107").Value = 917.6667
$ws.Range("M107").Value = 1508.1154
$ws.Range("N107").Value = -4757.6667

# Row 113
$ws.Range("H113").Value = 2562.611
$ws.Range("I113").Value = 1670
$ws.Range("K113").Value = 1670
$ws.Range("M113").Value = 1584

### Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")

# Row 45
$ws.Range("H45").Value = 1899.9286
$ws.Range("I45").Value = 1872.7273
$ws.Range("J45").Value = 1999.6666
$ws.Range("K45").Value = 1872.7273
$ws.Range("L45").Value = 1999.6666
$ws.Range("M45").Value = -1495.7273
$ws.Range("N45").Value = -2753.6666

# Row 88
$ws.Range("H88").Value = 12201.2
$ws.Range("I88").Value = 18002
$ws.Range("J88").Value = 3500
$ws.Range("K88").Value = 18002
$ws.Range("L88").Value = 3500
$ws.Range("M88").Value = -17596
$ws.Range("N88").Value = -4312

# Row 91
$ws.Range("H91").Value = 12201.2
$ws.Range("I91").Value = 18002
$ws.Range("J91").Value = 3500
$ws.Range("K91").Value = 18002
$ws.Range("L91").Value = 3500
$ws.Range("M91").Value = -16598
$ws.Range("N91").Value = -6308

### Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")

# Row 86
$ws.Range("H86").Value = 1876.5
$ws.Range("I86").Value = 1862.8
$ws.Range("J86").Value = 1967.8334
$ws.Range("K86").Value = 1862.8
$ws.Range("L86").Value = 1967.8334
$ws.Range("M86").Value = -739.8
$ws.Range("N86").Value = -4213.8334

# Row 89
$ws.Range("H89").Value = 1876.5
$ws.Range("I89").Value = 1862.8
$ws.Range("J89").Value = 1967.8334
$ws.Range("K89").Value = 9314
$ws.Range("L89").Value = 9839.166999999999
$ws.Range("M89").Value = -3698
$ws.Range("N89").Value = -21071.167

### Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")

# Row 58
$ws.Range("H58").Value = 2070319.2
$ws.Range("I58").Value = 2842510.8
$ws.Range("J58").Value = 11141.833
$ws.Range("K58").Value = 2842510.8
$ws.Range("L58").Value = 11141.833
$ws.Range("M58").Value = -2842307.8
$ws.Range("N58").Value = -11547.833

# Row 62
$ws.Range("H62").Value = 4000
$ws.Range("I62").Value = 4333.3335
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 4333.3335
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -3709.3335
$ws.Range("N62").Value = -4248

# Row 65
$ws.Range("H65").Value = 4000
$ws.Range("I65").Value = 4333.3335
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 21666.6675
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -18546.6675
$ws.Range("N65").Value = -21240

# Row 86
$ws.Range("H86").Value = 1788.4517
$ws.Range("I86").Value = 1755.5385
$ws.Range("K86").Value = 1755.5385
$ws.Range("M86").Value = -632.5385000000001

# Row 89
$ws.Range("H89").Value = 1788.4517
$ws.Range("I89").Value = 1755.5385
$ws.Range("K89").Value = 8777.692500000001
$ws.Range("M89").Value = -3161.692500000001

# Row 136
$ws.Range("H136").Value = 2070319.2
$ws.Range("I136").Value = 2842510.8
$ws.Range("J136").Value = 11141.833
$ws.Range("K136").Value = 8527532.399999999
$ws.Range("L136").Value = 33425.499
$ws.Range("M136").Value = -8524982.399999999
$ws.Range("N136").Value = -38525.499

### Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")

# Row 87
$ws.Range("H87").Value = 7722.7646
$ws.Range("I87").Value = 974.75
$ws.Range("J87").Value = 9799.076999999999
$ws.Range("K87").Value = 2924.25
$ws.Range("L87").Value = 29397.231
$ws.Range("M87").Value = -1676.25
$ws.Range("N87").Value = -31893.231

# Row 90
$ws.Range("H90").Value = 7722.7646
$ws.Range("I90").Value = 974.75
$ws.Range("J90").Value = 9799.076999999999
$ws.Range("K90").Value = 8772.75
$ws.Range("L90").Value = 88191.693
$ws.Range("M90").Value = -2532.75
$ws.Range("N90").Value = -100671.693

# Row 92
$ws.Range("H92").Value = 262.5
$ws.Range("I92").Value = 262.5
$ws.Range("K92").Value = 787.5
$ws.Range("M92").Value = 460.5

# Row 98
$ws.Range("H98").Value = 434.96295
$ws.Range("J98").Value = 918.8
$ws.Range("L98").Value = 2756.4
$ws.Range("N98").Value = -5752.4

# Row 120
$ws.Range("H120").Value = 8455.556
$ws.Range("I120").Value = 9733.333000000001
$ws.Range("J120").Value = 7816.6665
$ws.Range("K120").Value = 29199.999
$ws.Range("L120").Value = 23449.9995
$ws.Range("M120").Value = -24361.999
$ws.Range("N120").Value = -33125.99950000001

# Row 131
$ws.Range("H131").Value = 1440.1818
$ws.Range("J131").Value = 1092
$ws.Range("L131").Value = 3276
$ws.Range("N131").Value = -13356

# Row 132
$ws.Range("H132").Value = 1688.9
$ws.Range("I132").Value = 1889.5
$ws.Range("J132").Value = 1488.3
$ws.Range("K132").Value = 17005.5
$ws.Range("L132").Value = 13394.7
$ws.Range("M132").Value = -14475.5
$ws.Range("N132").Value = -18454.7

### Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")

# Row 51
$ws.Range("H51").Value = 24000
$ws.Range("J51").Value = 24000
$ws.Range("L51").Value = 24000
$ws.Range("N51").Value = -25018

# Row 68
$ws.Range("H68").Value = 40295
$ws.Range("J68").Value = 40295
$ws.Range("L68").Value = 40295
$ws.Range("N68").Value = -41917

# Row 70
$ws.Range("H70").Value = 5610.5713
$ws.Range("I70").Value = 5454.8
$ws.Range("J70").Value = 6000
$ws.Range("K70").Value = 5454.8
$ws.Range("L70").Value = 6000
$ws.Range("M70").Value = -5184.8
$ws.Range("N70").Value = -6540

# Row 71
$ws.Range("H71").Value = 40295
$ws.Range("J71").Value = 40295
$ws.Range("L71").Value = 120885
$ws.Range("N71").Value = -128997

# Row 73
$ws.Range("H73").Value = 5610.5713
$ws.Range("I73").Value = 5454.8
$ws.Range("J73").Value = 6000
$ws.Range("K73").Value = 5454.8
$ws.Range("L73").Value = 6000
$ws.Range("M73").Value = -4518.8
$ws.Range("N73").Value = -7872

# Row 80
$ws.Range("H80").Value = 11333.333
$ws.Range("I80").Value = 26500
$ws.Range("J80").Value = 7000
$ws.Range("K80").Value = 26500
$ws.Range("L80").Value = 7000
$ws.Range("M80").Value = -25502
$ws.Range("N80").Value = -8996

# Row 83
$ws.Range("H83").Value = 11333.333
$ws.Range("I83").Value = 26500
$ws.Range("J83").Value = 7000
$ws.Range("K83").Value = 132500
$ws.Range("L83").Value = 35000
$ws.Range("M83").Value = -127508
$ws.Range("N83").Value = -44984

### Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")

# Row 61
$ws.Range("H61").Value = 2550752.2
$ws.Range("I61").Value = 67668
$ws.Range("K61").Value = 67668
$ws.Range("M61").Value = -67466

# Row 100
$ws.Range("H100").Value = 4833.3335
$ws.Range("I100").Value = 4444.4443
$ws.Range("J100").Value = 6000
$ws.Range("K100").Value = 4444.4443
$ws.Range("L100").Value = 6000
$ws.Range("M100").Value = -3903.4443
$ws.Range("N100").Value = -7082

# Row 113
$ws.Range("H113").Value = 2550752.2
$ws.Range("I113").Value = 67668
$ws.Range("K113").Value = 67668
$ws.Range("M113").Value = -65498

### Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")

# Row 70
$ws.Range("H70").Value = 36232.855
$ws.Range("J70").Value = 36232.855
$ws.Range("L70").Value = 36232.855
$ws.Range("N70").Value = -36862.855

# Row 73
$ws.Range("H73").Value = 36232.855
$ws.Range("J73").Value = 36232.855
$ws.Range("L73").Value = 36232.855
$ws.Range("N73").Value = -38416.855

# Row 80
$ws.Range("H80").Value = 40301
$ws.Range("J80").Value = 40301
$ws.Range("L80").Value = 40301
$ws.Range("N80").Value = -42297

# Row 83
$ws.Range("H83").Value = 40301
$ws.Range("J83").Value = 40301
$ws.Range("L83").Value = 120903
$ws.Range("N83").Value = -130887

